$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Range("H80").Value = 4705.722
$ws.Range("J80").Value = 6766.9165
$ws.Range("L80").Value = 20300.7495
$ws.Range("N80").Value = -22296.7495

# Row 83
$ws.Range("H83").Value = 4705.722
$ws.Range("J83").Value = 6766.9165
$ws.Range("L83").Value = 60902.2485
$ws.Range("N83").Value = -70886.2485

# Row 86
$ws.Range("H86").Value = 5845.04
$ws.Range("I86").Value = 21262.4
$ws.Range("J86").Value = 1990.7
$ws.Range("K86").Value = 21262.4
$ws.Range("L86").Value = 1990.7
$ws.Range("M86").Value = -20139.4
$ws.Range("N86").Value = -4236.7

# Row 89
$ws.Range("H89").Value = 5845.04
$ws.Range("I89").Value = 21262.4
$ws.Range("J89").Value = 1990.7
$ws.Range("K89").Value = 106312
$ws.Range("L89").Value = 9953.5
$ws.Range("M89").Value = -100696
$ws.Range("N89").Value = -21185.5

# Row 98
$ws.Range("H98").Value = 1886.5428
$ws.Range("I98").Value = 1380.9667
$ws.Range("J98").Value = 4920
$ws.Range("K98").Value = 1380.9667
$ws.Range("L98").Value = 4920
$ws.Range("M98").Value = 117.0333000000001
$ws.Range("N98").Value = -7916

# Row 122
$ws.Range("H122").Value = 1886.5428
$ws.Range("I122").Value = 1380.9667
$ws.Range("J122").Value = 4920
$ws.Range("K122").Value = 4142.9001
$ws.Range("L122").Value = 14760
$ws.Range("M122").Value = -1692.9001
$ws.Range("N122").Value = -19660

$ws = $wb.Worksheets.Item("ARM")
# Row 52
$ws.Range("H52").Value = 39779.41
$ws.Range("J52").Value = 39779.41
$ws.Range("L52").Value = 39779.41
$ws.Range("N52").Value = -40415.41

# Row 61
$ws.Range("H61").Value = 1964.3158
$ws.Range("I61").Value = 1038.125
$ws.Range("J61").Value = 6904
$ws.Range("K61").Value = 1038.125
$ws.Range("L61").Value = 6904
$ws.Range("M61").Value = -826.125
$ws.Range("N61").Value = -7328

# Row 64
$ws.Range("H64").Value = 28860
$ws.Range("I64").Value = 26000
$ws.Range("J64").Value = 29177.777
$ws.Range("K64").Value = 26000
$ws.Range("L64").Value = 29177.777
$ws.Range("M64").Value = -25752
$ws.Range("N64").Value = -29673.777

# Row 67
$ws.Range("H67").Value = 28860
$ws.Range("I67").Value = 26000
$ws.Range("J67").Value = 29177.777
$ws.Range("K67").Value = 26000
$ws.Range("L67").Value = 29177.777
$ws.Range("M67").Value = -25142
$ws.Range("N67").Value = -30893.777

# Row 122
$ws.Range("H122").Value = 1950.1428
$ws.Range("I122").Value = 1061.8
$ws.Range("J122").Value = 2443.6667
$ws.Range("K122").Value = 3185.4
$ws.Range("L122").Value = 7331.000100000001
$ws.Range("M122").Value = -735.3999999999996
$ws.Range("N122").Value = -12231.0001

# Row 132
$ws.Range("H132").Value = 20836820
$ws.Range("I132").Value = 23812800
$ws.Range("J132").Value = 4949.6665
$ws.Range("K132").Value = 71438400
$ws.Range("L132").Value = 14848.9995
$ws.Range("M132").Value = -71435870
$ws.Range("N132").Value = -19908.9995

# Row 136
$ws.Range("H136").Value = 1964.3158
$ws.Range("I136").Value = 1038.125
$ws.Range("J136").Value = 6904
$ws.Range("K136").Value = 3114.375
$ws.Range("L136").Value = 20712
$ws.Range("M136").Value = -564.375
$ws.Range("N136").Value = -25812

# Row 139
$ws.Range("H139").Value = 85000
$ws.Range("J139").Value = 85000
$ws.Range("L139").Value = 85000
$ws.Range("N139").Value = -95280

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1592.7037
$ws.Range("I20").Value = 1021.1875
$ws.Range("K20").Value = 1021.1875
$ws.Range("M20").Value = -774.1875

# Row 134
$ws.Range("H134").Value = 2085.162
$ws.Range("I134").Value = 1438.3667
$ws.Range("K134").Value = 4315.1001
$ws.Range("M134").Value = -1780.1001

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1890217.2
$ws.Range("I31").Value = 2501344.5
$ws.Range("J31").Value = 9825.308000000001
$ws.Range("K31").Value = 2501344.5
$ws.Range("L31").Value = 9825.308000000001
$ws.Range("M31").Value = -2501049.5
$ws.Range("N31").Value = -10415.308

# Row 34
$ws.Range("H34").Value = 1890217.2
$ws.Range("I34").Value = 2501344.5
$ws.Range("J34").Value = 9825.308000000001
$ws.Range("K34").Value = 2501344.5
$ws.Range("L34").Value = 9825.308000000001
$ws.Range("M34").Value = -2501142.5
$ws.Range("N34").Value = -10229.308

$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 1512.125
$ws.Range("I107").Value = 1024.25
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 3072.75
$ws.Range("L107").Value = 6000
$ws.Range("M107").Value = -1152.75
$ws.Range("N107").Value = -9840

# Row 118
$ws.Range("H118").Value = 3937.375
$ws.Range("I118").Value = 719.4
$ws.Range("J118").Value = 4784.2104
$ws.Range("K118").Value = 2158.2
$ws.Range("L118").Value = 14352.6312
$ws.Range("M118").Value = -915.1999999999998
$ws.Range("N118").Value = -16838.6312

# Row 122
$ws.Range("H122").Value = 1189.1333
$ws.Range("I122").Value = 486.66666
$ws.Range("J122").Value = 1364.75
$ws.Range("K122").Value = 4379.99994
$ws.Range("L122").Value = 12282.75
$ws.Range("M122").Value = -1929.99994
$ws.Range("N122").Value = -17182.75

$ws = $wb.Worksheets.Item("LTW")
# Row 64
$ws.Range("H64").Value = 14136
$ws.Range("I64").Value = 14136
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 14136
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -13911

# Row 67
$ws.Range("H67").Value = 14136
$ws.Range("I67").Value = 14136
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 14136
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -13356

# Row 68
$ws.Range("H68").Value = 1972.1111
$ws.Range("I68").Value = 978.4286
$ws.Range("J68").Value = 5450
$ws.Range("K68").Value = 978.4286
$ws.Range("L68").Value = 5450
$ws.Range("M68").Value = -229.4286
$ws.Range("N68").Value = -6948

# Row 71
$ws.Range("H71").Value = 1972.1111
$ws.Range("I71").Value = 978.4286
$ws.Range("J71").Value = 5450
$ws.Range("K71").Value = 4892.143
$ws.Range("L71").Value = 27250
$ws.Range("M71").Value = -1148.143
$ws.Range("N71").Value = -34738

# Row 93
$ws.Range("H93").Value = 2232.75
$ws.Range("I93").Value = 1448.625
$ws.Range("J93").Value = 3801
$ws.Range("K93").Value = 1448.625
$ws.Range("L93").Value = 3801
$ws.Range("M93").Value = -200.625
$ws.Range("N93").Value = -6297

# Row 106
$ws.Range("H106").Value = 24000
$ws.Range("J106").Value = 24000
$ws.Range("L106").Value = 24000
$ws.Range("N106").Value = -26524

# Row 122
$ws.Range("H122").Value = 2802.6584
$ws.Range("I122").Value = 2682.9714
$ws.Range("J122").Value = 3500.8333
$ws.Range("K122").Value = 8048.914199999999
$ws.Range("L122").Value = 10502.4999
$ws.Range("M122").Value = -5598.914199999999
$ws.Range("N122").Value = -15402.4999

# Row 140
$ws.Range("H140").Value = 30000
$ws.Range("J140").Value = 30000
$ws.Range("L140").Value = 30000
$ws.Range("N140").Value = -40360

# Row 141
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").ClearContents()
$ws.Range("N141").Value = 0

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 308390.56
$ws.Range("I132").Value = 477759.8
$ws.Range("J132").Value = 11994.417
$ws.Range("K132").Value = 1433279.4
$ws.Range("L132").Value = 35983.251
$ws.Range("M132").Value = -1430749.4
$ws.Range("N132").Value = -41043.251
